# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list ... with GitHub Actions" (Sun Dec 17 14:46:39 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text such as "41.828.86" or "0.624". Excel would
# silently reinterpret digit-and-dot strings as numbers, so force the cell
# to text first (and reset the style afterwards so no extra formatting is
# left behind) before writing the value - this mirrors the original sheet,
# where every price is stored as plain text.
function Set-TextValue([string]$cellRef, [string]$val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '41.828.86'
$ws.Range("E2").Value = '  -1.41%  '
Set-TextValue "D3" '2.214.96'
$ws.Range("E3").Value = '  -1.36%  '
$ws.Range("E4").Value = '  +0.10%  '
Set-TextValue "D5" '241.22'
$ws.Range("E5").Value = '  -1.99%  '
Set-TextValue "D6" '0.624'
$ws.Range("E6").Value = '  +0.27%  '
Set-TextValue "D7" '72.57'
$ws.Range("E7").Value = '  -4.19%  '
$ws.Range("E8").Value = '  +0.16%  '
Set-TextValue "D9" '0.597'
$ws.Range("E9").Value = '  -3.58%  '
Set-TextValue "D10" '41.64'
$ws.Range("E10").Value = '  -4.31%  '
Set-TextValue "D11" '0.0944'
$ws.Range("E11").Value = '  -0.42%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D12" '6.96'
$ws.Range("E12").Value = '  -3.92%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue "D13" '0.103'
$ws.Range("E13").Value = '  +0.50%  '
Set-TextValue "D14" '2.545.18'
$ws.Range("E14").Value = '  -1.37%  '
Set-TextValue "D15" '14.18'
$ws.Range("E15").Value = '  -2.69%  '
Set-TextValue "D16" '0.832'
$ws.Range("E16").Value = '  -2.87%  '
Set-TextValue "D17" '2.196.66'
$ws.Range("E17").Value = '  -2.43%  '
Set-TextValue "D18" '41.689.47'
$ws.Range("E18").Value = '  -1.17%  '
Set-TextValue "D19" '0.0000105'
$ws.Range("E19").Value = '  +2.54%  '
Set-TextValue "D20" '72.31'
$ws.Range("E20").Value = '  +0.24%  '
Set-TextValue "D21" '6.14'
$ws.Range("E21").Value = '  -0.47%  '
Set-TextValue "D22" '11.13'
$ws.Range("E22").Value = '  +20.63%  '
Set-TextValue "D23" '228.81'
$ws.Range("E23").Value = '  -0.54%  '
Set-TextValue "D24" '2.04'
$ws.Range("E24").Value = '  -8.84%  '
$ws.Range("E25").Value = '  +0.25%  '
Set-TextValue "D26" '11.34'
$ws.Range("E26").Value = '  -1.05%  '
Set-TextValue "D27" '3.63'
$ws.Range("E27").Value = '  +0.11%  '
Set-TextValue "D28" '2.27'
$ws.Range("E28").Value = '  -1.69%  '
$ws.Range("E29").Value = '  -0.63%  '
Set-TextValue "D30" '167.03'
$ws.Range("E30").Value = '  -0.37%  '
Set-TextValue "D31" '20.38'
$ws.Range("E31").Value = '  -1.43%  '
Set-TextValue "D32" '5.60'
$ws.Range("E32").Value = '  +5.21%  '
Set-TextValue "D33" '0.0795'
$ws.Range("E33").Value = '  -3.56%  '
Set-TextValue "D34" '29.97'
$ws.Range("E34").Value = '  -3.02%  '
$ws.Range("E35").Value = '  -0.99%  '
Set-TextValue "D36" '0.107'
$ws.Range("E36").Value = '  -10.99%  '
Set-TextValue "D37" '4.21'
$ws.Range("E37").Value = '  -6.91%  '
Set-TextValue "D38" '0.0299'
$ws.Range("E38").Value = '  -5.06%  '
Set-TextValue "D39" '13.38'
$ws.Range("E39").Value = '  -4.21%  '
Set-TextValue "D40" '2.12'
$ws.Range("E40").Value = '  -2.88%  '
Set-TextValue "D41" '5.60'
$ws.Range("E41").Value = '  -3.39%  '
Set-TextValue "D42" '63.83'
$ws.Range("E42").Value = '  -0.50%  '
Set-TextValue "D43" '0.196'
$ws.Range("E43").Value = '  -2.61%  '
Set-TextValue "D44" '8.65'
$ws.Range("E44").Value = '  -1.65%  '
Set-TextValue "D45" '102.89'
$ws.Range("E45").Value = '  -4.65%  '
Set-TextValue "D46" '0.0998'
$ws.Range("E46").Value = '  -2.86%  '
Set-TextValue "D47" '2.33'
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("E48").Value = '  -1.98%  '
$ws.Range("E49").Value = '  -2.56%  '
Set-TextValue "D50" '2.68'
$ws.Range("E50").Value = '  -0.79%  '
Set-TextValue "D51" '2.416.52'
$ws.Range("E51").Value = '  -1.63%  '
